$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price figures that must remain plain text even
# though many of them look like numbers (e.g. "3.08", "0.999"). Force the
# cell to text format before assigning so Excel does not coerce the string
# into a numeric value, then restore the default "Normal" style so no stray
# formatting is left behind (matches the source workbook, which carries no
# explicit style on these cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.430.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.775.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0835"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.212.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.774.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.922"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.414.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.162"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +6.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("E35").Value = "  +5.94%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.113.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.902"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("E51").Value = "  +7.18%  "
